$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new value looks like a plain number,
# so Excel stores them as text (matching the source data) instead of converting to a number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '29.534.82'
$ws.Range("E2").Value = '  +1.03%  '

# Row 3
$ws.Range("D3").Value = '1.877.35'
$ws.Range("E3").Value = '  +0.75%  '

# Row 4
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").Value = '0.7219'
$ws.Range("E5").Value = '  +1.78%  '

# Row 6
$ws.Range("D6").Value = '239.62'
$ws.Range("E6").Value = '  +0.78%  '

# Row 7
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.13%  '

# Row 8
$ws.Range("D8").Value = '0.07867'
$ws.Range("E8").Value = '  -4.12%  '

# Row 9
$ws.Range("D9").Value = '0.3090'
$ws.Range("E9").Value = '  +1.73%  '

# Row 10
$ws.Range("E10").Value = '  +8.29%  '

# Row 11
$ws.Range("D11").Value = '0.08221'
$ws.Range("E11").Value = '  +0.61%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.891.70'
$ws.Range("E12").Value = '  +2.97%  '

# Row 13
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.7263'
$ws.Range("E13").Value = '  +2.50%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.260'
$ws.Range("E14").Value = '  +1.76%  '

# Row 15
$ws.Range("D15").Value = '90.30'
$ws.Range("E15").Value = '  +1.13%  '

# Row 16
$ws.Range("D16").Value = '29.497.70'
$ws.Range("E16").Value = '  +0.93%  '

# Row 17
$ws.Range("D17").Value = '5.863'
$ws.Range("E17").Value = '  +1.33%  '

# Row 18
$ws.Range("D18").Value = '0.000007875'
$ws.Range("E18").Value = '  -0.31%  '

# Row 19
$ws.Range("E19").Value = '  +2.31%  '

# Row 20
$ws.Range("D20").Value = '13.37'
$ws.Range("E20").Value = '  +0.17%  '

# Row 21
$ws.Range("D21").Value = '2.137.91'
$ws.Range("E21").Value = '  +1.63%  '

# Row 22
$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  +0.11%  '

# Row 23
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.21%  '

# Row 24
$ws.Range("D24").Value = '7.761'
$ws.Range("E24").Value = '  +4.84%  '

# Row 25
$ws.Range("D25").Value = '0.1591'
$ws.Range("E25").Value = '  +10.16%  '

# Row 26
$ws.Range("D26").Value = '162.96'
$ws.Range("E26").Value = '  +0.47%  '

# Row 27
$ws.Range("D27").Value = '8.995'
$ws.Range("E27").Value = '  +0.33%  '

# Row 28
$ws.Range("D28").Value = '18.39'
$ws.Range("E28").Value = '  +1.69%  '

# Row 29
$ws.Range("D29").Value = '1.947'
$ws.Range("E29").Value = '  -0.72%  '

# Row 30
$ws.Range("D30").Value = '1.371'
$ws.Range("E30").Value = '  -3.73%  '

# Row 31
$ws.Range("D31").Value = '1.486'
$ws.Range("E31").Value = '  +0.09%  '

# Row 32
$ws.Range("D32").Value = '4.360'
$ws.Range("E32").Value = '  -0.63%  '

# Row 33
$ws.Range("D33").Value = '4.092'
$ws.Range("E33").Value = '  +0.96%  '

# Row 34
$ws.Range("D34").Value = '0.05262'
$ws.Range("E34").Value = '  +1.03%  '

# Row 35
$ws.Range("E35").Value = '  +2.35%  '

# Row 36
$ws.Range("D36").Value = '0.7189'
$ws.Range("E36").Value = '  +1.66%  '

# Row 37
$ws.Range("D37").Value = '0.9998'
$ws.Range("E37").Value = '  +0.09%  '

# Row 38
$ws.Range("D38").Value = '2.671'
$ws.Range("E38").Value = '  +0.06%  '

# Row 39
$ws.Range("D39").Value = '0.01865'
$ws.Range("E39").Value = '  +1.01%  '

# Row 40
$ws.Range("D40").Value = '2.716'
$ws.Range("E40").Value = '  -0.38%  '

# Row 41
$ws.Range("D41").Value = '1.185.03'
$ws.Range("E41").Value = '  +3.70%  '

# Row 42
$ws.Range("D42").Value = '0.9011'
$ws.Range("E42").Value = '  -2.54%  '

# Row 43
$ws.Range("D43").Value = '6.007'
$ws.Range("E43").Value = '  +2.44%  '

# Row 44
$ws.Range("D44").Value = '0.4329'
$ws.Range("E44").Value = '  +1.25%  '

# Row 45
$ws.Range("D45").Value = '71.98'
$ws.Range("E45").Value = '  +2.59%  '

# Row 46
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.17%  '

# Row 47
$ws.Range("D47").Value = '102.72'
$ws.Range("E47").Value = '  +0.02%  '

# Row 48
$ws.Range("D48").Value = '0.5363'
$ws.Range("E48").Value = '  -0.85%  '

# Row 49
$ws.Range("D49").Value = '1.777'
$ws.Range("E49").Value = '  +0.04%  '

# Row 50
$ws.Range("D50").Value = '9.251'
$ws.Range("E50").Value = '  +0.50%  '

# Row 51
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").Value = '2.880'
$ws.Range("E51").Value = '  +5.02%  '
